# UInfo.xlsx - "segundo xlsx y mejora del programa"
# Rework the interaction-count sheet: extend the step table from 5 steps to a
# generic up-to-10-steps table, add a new "h" column/header, rebuild the
# weighting/selection rows, and move the explanatory note + the little
# X/Y legend to the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- grab the formats we still need before their old homes are repurposed -
# the boxed note (D6, merged D6:H6) moves down to K8 (merged K8:O8); the
# light "fill" cells next to it (E6:H6) supply L8:O8's format.
$ws.Range("D6").Copy()
$ws.Range("K8").PasteSpecial(-4122) | Out-Null
$ws.Range("E6").Copy()
$ws.Range("L8:O8").PasteSpecial(-4122) | Out-Null

# the bordered "step value" look used by D4:H4 / D5:H5 extends out to M
$ws.Range("H4").Copy()
$ws.Range("I4:M4").PasteSpecial(-4122) | Out-Null
$ws.Range("H5").Copy()
$ws.Range("I5:M5").PasteSpecial(-4122) | Out-Null
$ws.Range("C5").Copy()
$ws.Range("C6").PasteSpecial(-4122) | Out-Null

# --- remove the things that are going away -----------------------------
$ws.Range("D6:H6").UnMerge()
$ws.Range("H15:I15").Clear()

# put D6:M6 back to the sheet's plain default look (no border/fill)
$ws.Range("C3").Copy()
$ws.Range("D6:M6").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# --- header row (row 2) --------------------------------------------------
$ws.Range("C2").Value = "n"
$ws.Range("D2").Value = "steps"
$ws.Range("E2").Value = "b"
$ws.Range("F2").Value = "h"

# --- parameter row (row 3) -----------------------------------------------
$ws.Range("C3").Value = 4
$ws.Range("D3").Value = 0.25
$ws.Range("E3").Value = 2
$ws.Range("F3").Formula = '=(E3-C4)/2'

# --- "X" step-value row (row 4) ------------------------------------------
$ws.Range("B4").Value = "X"
$ws.Range("C4").Value = 1
$ws.Range("D4").Formula = '=IF(C3>=1,C4+D3,"")'
$ws.Range("E4").Formula = '=IF(C3>=2,D4+D3,"")'
$ws.Range("F4").Formula = '=IF(C3>=3,E4+D3,"")'
$ws.Range("G4").Formula = '=IF(C3>=4,F4+D3,"")'
$ws.Range("H4").Formula = '=IF(C3>=5,G4+D3,"")'
$ws.Range("I4").Formula = '=IF(C3>=6,H4+D3,"")'
$ws.Range("J4").Formula = '=IF(C3>=7,I4+D3,"")'
$ws.Range("K4").Formula = '=IF(C3>=8,J4+D3,"")'
$ws.Range("L4").Formula = '=IF(C3>=9,K4+D3,"")'
$ws.Range("M4").Formula = '=IF(C3>=10,L4+D3,"")'

# --- "Y" weighting row (row 5) -------------------------------------------
$ws.Range("B5").Value = "Y"
$ws.Range("C5").Formula = '=(((2*C4)+1)*(C4-2)/C4)'
$ws.Range("D5").Formula = '=IF(D4="","",(((2*D4)+1)*((D4-2))/D4))'
$ws.Range("E5").Formula = '=IF(E4="","",(((2*E4)+1)*(E4-2)/E4))'
$ws.Range("F5").Formula = '=IF(F4="","",(((2*F4)+1)*(F4-2)/F4))'
$ws.Range("G5").Formula = '=IF(G4="","",(((2*G4)+1)*(G4-2)/G4))'
$ws.Range("H5:K5").Formula = '=IF(H4="","",1)'
$ws.Range("L5").Formula = '=IF(L4="","",1)'
$ws.Range("M5").Formula = '=IF(M4="","",1)'

# --- selector row (row 6) -------------------------------------------------
$ws.Range("D6").Formula = '=IF(OR(E5="",C3=1),IF(AND(E5="",C3=1),1,0),IF(C3=2,4,1))'
$ws.Range("E6").Formula = '=IF(OR(F5="",C3=2),IF(AND(F5="",C3=2),1,0),IF(C3=4,4,1))'
$ws.Range("F6").Formula = '=IF(OR(G5="",C3=3),IF(AND(G5="",C3=3),1,0),IF(C3=6,4,1))'
$ws.Range("G6").Formula = '=IF(OR(H5="",C3=4),IF(AND(H5="",C3=4),1,0),IF(C3=8,4,1))'
$ws.Range("H6").Formula = '=IF(OR(I5="",C3=5),IF(AND(I5="",C3=5),1,0),IF(C3=10,4,1))'
$ws.Range("I6").Formula = '=IF(OR(J5="",C3=6),IF(AND(J5="",C3=6),1,0),1)'
$ws.Range("J6").Formula = '=IF(OR(K5="",C3=7),IF(AND(K5="",C3=7),1,0),1)'
$ws.Range("K6").Formula = '=IF(OR(L5="",C3=8),IF(AND(L5="",C3=8),1,0),1)'
$ws.Range("L6").Formula = '=IF(OR(M5="",C3=9),IF(AND(M5="",C3=9),1,0),1)'
$ws.Range("M6").Value = 1

# --- final result + note (row 8) -----------------------------------------
$ws.Range("E8").Formula = '=(F3/3)*(C5+IF(D4="",0,IF(C3=2,(D5*D6),0))+IF(E4="",0,IF(C3=4,(E5*E6),IF(C3=2,E5)))+IF(F4="",0,IF(C3=6,(F5*F6),0)+IF(G4="",0,IF(C3=8,(G5*G6),IF(C3=4,G5)+IF(H4="",0,IF(C3=10,(H5*H6),0)+IF(I4="",0,IF(C3=6,(I5*I6)))+IF(J4="",0,IF(C3=8,(J5*J6)))+IF(K4="",0,IF(C3=10,(K5*K6)))+IF(L4="",0,(L5*L6))+IF(M4="",0,(M5*M6)))))))'
$ws.Range("K8").Value = "Estas serian las celdas que cuentan como interacciones par"
$ws.Range("K8:O8").Merge() | Out-Null

# --- helper / substitution test cells (rows 11-12) ------------------------
$ws.Range("H11").Formula = '=X2+4'
$ws.Range("H12").Formula = '=SUBSTITUTE(H11,"x",D4)'

# --- legend, moved to the bottom of the sheet (row 22) --------------------
$ws.Range("O22").Value = " "
$ws.Range("P22").Value = "   "

# --- column sizing to mirror the authored layout --------------------------
$ws.Columns.Item(3).ColumnWidth = 13.5
$ws.Columns.Item(4).ColumnWidth = 13.95
$ws.Columns.Item(8).ColumnWidth = 13.5

# --- selection cursor, as left by the author ------------------------------
$ws.Range("E3").Select() | Out-Null
